# Weekly fruit/vegetable price update: insert a new observation row for
# "Mora" (Mercado Mayorista Lo Valledor de Santiago) at row 41, pushing the
# existing rows 41-86 down to 42-87.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41 (shifts rows 41-86 -> 42-87, carrying
# along their values/formatting, same as Excel's Insert Row behaviour).
$ws.Rows.Item(41).Insert()

# Populate the newly-inserted row 41 with the new weekly record. Most
# columns mirror the (now-shifted) row that used to be row 41 / is now
# row 42; only the date, min/max/avg price, origin and $/Kg columns differ.
$ws.Cells.Item(41, 1).Value = 6
$ws.Cells.Item(41, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(41, 3).Value = 'Metropolitana'
$ws.Cells.Item(41, 4).Value = 44894
$ws.Cells.Item(41, 5).Value = 13
$ws.Cells.Item(41, 6).Value = 'Fruta'
$ws.Cells.Item(41, 7).Value = 100101
$ws.Cells.Item(41, 8).Value = 'Berries'
$ws.Cells.Item(41, 9).Value = 100101008
$ws.Cells.Item(41, 10).Value = 'Mora'
$ws.Cells.Item(41, 11).Value = 'Sin especificar'
$ws.Cells.Item(41, 12).Value = 'Primera'
$ws.Cells.Item(41, 13).Value = 250
$ws.Cells.Item(41, 14).Value = 7000
$ws.Cells.Item(41, 15).Value = 7000
$ws.Cells.Item(41, 16).Value = 7000
$ws.Cells.Item(41, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(41, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(41, 19).Value = 3500
$ws.Cells.Item(41, 20).Value = 2
